$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 27 - pushes existing rows 27..98 down to 28..99,
# matching the row-level shift seen throughout the diff (e.g. old D27=44244
# reappears as the new D28, old D97=44580 reappears as the new D98, etc.)
$ws.Rows(27).Insert()

# Populate the freshly inserted row 27 with the new record.
$ws.Range("A27").Value = 2
$ws.Range("B27").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C27").Value = "Coquimbo"
$ws.Range("D27").Value = 44910
$ws.Range("E27").Value = 4
$ws.Range("F27").Value = 100112030
$ws.Range("G27").Value = "Poroto granado"
$ws.Range("H27").Value = "Sin especificar"
$ws.Range("I27").Value = "Primera"
$ws.Range("J27").Value = 500
$ws.Range("K27").Value = 28000
$ws.Range("L27").Value = 30000
$ws.Range("M27").Value = 29000
$ws.Range("N27").Value = "$/malla 25 kilos"
$ws.Range("O27").Value = "Provincia de Limarí"
$ws.Range("P27").Value = 1160
$ws.Range("Q27").Value = 25
$ws.Range("R27").Value = "Hortaliza"
